$wb = $excel.ActiveWorkbook

# --- Sheet 1: "device" -> clear all header content (id/name/description removed) ---
$deviceWs = $wb.Worksheets.Item("device")
$deviceWs.Cells.Clear()

# --- Sheet 2: "NamedThing" -> renamed to "catalysisSample" with new header row
#              + a list data-validation on column E (vital_status) ---
$catalysisWs = $wb.Worksheets.Item("NamedThing")
$catalysisWs.Cells.Clear()
$catalysisWs.Cells.Validation.Delete()
$catalysisWs.Range("A1").Value = "sample_environment"
$catalysisWs.Range("B1").Value = "primary_email"
$catalysisWs.Range("C1").Value = "birth_date"
$catalysisWs.Range("D1").Value = "age_in_years"
$catalysisWs.Range("E1").Value = "vital_status"
$dv = $catalysisWs.Range("E2:E1048576").Validation
$dv.Add(3, 1, 1, '"ALIVE,DEAD,UNKNOWN"')
$dv.InCellDropdown = $true
$dv.ShowInput = $false
$dv.ShowError = $false
$catalysisWs.Name = "catalysisSample"

# --- Sheet 3: "Sample" -> renamed to "sample" with trimmed header row
#              (id/name/description dropped, the vital_status validation dropped) ---
$sampleWs = $wb.Worksheets.Item("Sample")
$sampleWs.Cells.Clear()
$sampleWs.Cells.Validation.Delete()
$sampleWs.Range("A1").Value = "primary_email"
$sampleWs.Range("B1").Value = "birth_date"
$sampleWs.Range("C1").Value = "age_in_years"
$sampleWs.Range("D1").Value = "vital_status"
$sampleWs.Name = "sample"

# --- Sheet 4: "SampleCollection" -> renamed to "sampleCollection", content unchanged ---
$sampleCollectionWs = $wb.Worksheets.Item("SampleCollection")
$sampleCollectionWs.Name = "sampleCollection"

# --- Sheet 5 (new): "namedThing", empty sheet appended at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$namedThingWs = $wb.Worksheets.Add($null, $lastSheet)
$namedThingWs.Name = "namedThing"
